$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1895910780669145
    "C2" = 0.5539033457249071
    "J2" = 0.007434944237918215
    "P2" = 0.1338289962825279
    "S2" = 0.1152416356877323
    "C3" = 0.03821656050955414
    "J3" = 0.06369426751592357
    "P3" = 0.6687898089171974
    "S3" = 0.2292993630573248
    "J4" = 0.04166666666666666
    "P4" = 0.7083333333333334
    "S4" = 0.25
    "B6" = 0.06367041198501873
    "D6" = 0.003745318352059925
    "F6" = 0.0449438202247191
    "J6" = 0.2696629213483146
    "O6" = 0.0299625468164794
    "Q6" = 0.1947565543071161
    "R6" = 0.06741573033707865
    "S6" = 0.3258426966292135
    "B7" = 0.07373271889400922
    "D7" = 0.01382488479262673
    "F7" = 0.05529953917050692
    "J7" = 0.152073732718894
    "O7" = 0.01382488479262673
    "Q7" = 0.1889400921658986
    "R7" = 0.1013824884792627
    "S7" = 0.4009216589861751
    "B8" = 0.0735930735930736
    "D8" = 0.02597402597402598
    "F8" = 0.04761904761904762
    "J8" = 0.1147186147186147
    "O8" = 0.01948051948051948
    "Q8" = 0.1861471861471861
    "R8" = 0.09956709956709957
    "S8" = 0.4329004329004329
    "B9" = 0.07392996108949416
    "D9" = 0.007782101167315175
    "F9" = 0.07392996108949416
    "J9" = 0.1050583657587549
    "O9" = 0.03891050583657588
    "Q9" = 0.1750972762645914
    "R9" = 0.07392996108949416
    "S9" = 0.4513618677042802
    "B10" = 0.1005372217958557
    "D10" = 0.02455871066768995
    "E10" = 0.001534919416730622
    "F10" = 0.08135072908672294
    "J10" = 0.1105141980046048
    "O10" = 0.01381427475057559
    "Q10" = 0.2110514198004605
    "R10" = 0.07214121258633922
    "S10" = 0.3844973138910207
    "G11" = 0.09333333333333334
    "J11" = 0.1166666666666667
    "K11" = 0.1333333333333333
    "L11" = 0.6566666666666666
    "G12" = 0.7537688442211056
    "J12" = 0.1909547738693467
    "K12" = 0.005025125628140704
    "L12" = 0.01507537688442211
    "S12" = 0.03517587939698492
    "F13" = 0.01612903225806452
    "G13" = 0.7096774193548387
    "J13" = 0.2419354838709677
    "S13" = 0.03225806451612903
    "F15" = 0.01276595744680851
    "H15" = 0.1574468085106383
    "I15" = 0.0851063829787234
    "J15" = 0.3404255319148936
    "K15" = 0.04680851063829787
    "M15" = 0.01276595744680851
    "N15" = 0.00425531914893617
    "O15" = 0.06808510638297872
    "S15" = 0.2723404255319149
    "F16" = 0.02923976608187134
    "H16" = 0.1695906432748538
    "I16" = 0.09941520467836257
    "J16" = 0.3859649122807017
    "K16" = 0.1228070175438596
    "M16" = 0.04093567251461988
    "O16" = 0.02339181286549707
    "S16" = 0.1286549707602339
    "F17" = 0.03024193548387097
    "H17" = 0.1834677419354839
    "I17" = 0.1129032258064516
    "J17" = 0.3568548387096774
    "K17" = 0.1108870967741935
    "M17" = 0.02016129032258064
    "O17" = 0.07862903225806452
    "S17" = 0.1068548387096774
    "F18" = 0.03045685279187817
    "H18" = 0.1624365482233502
    "I18" = 0.1015228426395939
    "J18" = 0.4060913705583756
    "K18" = 0.1116751269035533
    "M18" = 0.03045685279187817
    "N18" = 0.005076142131979695
    "O18" = 0.04568527918781726
    "S18" = 0.1065989847715736
    "F19" = 0.02046783625730994
    "H19" = 0.2010233918128655
    "I19" = 0.1067251461988304
    "J19" = 0.3567251461988304
    "K19" = 0.1067251461988304
    "M19" = 0.02704678362573099
    "N19" = 0.0007309941520467836
    "O19" = 0.07383040935672515
    "S19" = 0.1067251461988304
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

Write-Output ("Updated " + $values.Count + " cells")
